# Contest 19 RCB vs DC
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the results for Contest 19 (row 28, "RCB vs DC")
$ws.Range("E28").Value = 100
$ws.Range("H28").Value = 60
$ws.Range("K28").Value = 40
$ws.Range("N28").Value = 80
$ws.Range("Q28").Value = 0
$ws.Range("T28").Value = 20

# Insert a new blank contest row before row 38, shifting the summary block down
$ws.Rows("38:38").Insert()

# Add Contest 28 ("RCB vs KKR") into row 37
$ws.Range("A37").Value = 28
$ws.Range("B37").Value = 1
$ws.Range("C37").Value = "RCB vs KKR"

$ws.Range("U42").Select()
